# This script applies a reshuffle of the daily price records (rows 2-10)
# in the "Hortaliza, Terminal Hortofrutícola Agro Chillán - Cebollín" sheet.
# Columns A, B, C, E, F, G, H, I, R are identical for every data row, so only
# the per-record columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), N (Unidad de
# comercializacion), O (Origen), P (Precio $/Kg) and Q (Kg o Unidades) need
# to move between rows, which is exactly what is reflected in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values travel together as one "record" per row.
$cols = @("D", "J", "K", "L", "M", "N", "O", "P", "Q")

# Maps destination row -> source row (both refer to the ORIGINAL/before state).
$rowMap = @{
    2  = 3
    3  = 10
    4  = 9
    5  = 2
    6  = 8
    7  = 5
    8  = 6
    9  = 4
    10 = 7
}

# Snapshot the original values for every relevant cell before overwriting
# anything, since several rows read from each other.
$snapshot = @{}
foreach ($r in 2..10) {
    foreach ($c in $cols) {
        $snapshot["$c$r"] = $ws.Range("$c$r").Value2
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $snapshot["$c$srcRow"]
    }
}
